$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.73"
$ws.Range("E2").Value = "'0.22%"
$ws.Range("D3").Value = "'41.06"
$ws.Range("E3").Value = "'-0.64%"
$ws.Range("D4").Value = "'5.227"
$ws.Range("E4").Value = "'1.86%"
$ws.Range("D5").Value = "'0.07667"
$ws.Range("E5").Value = "'0.68%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.618"
$ws.Range("E6").Value = "'-0.27%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9191"
$ws.Range("E7").Value = "'1.74%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.428"
$ws.Range("E8").Value = "'-1.93%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1241"
$ws.Range("E9").Value = "'13.40%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1839"
$ws.Range("E10").Value = "'3.97%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09175"
$ws.Range("E11").Value = "'0.18%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04251"
$ws.Range("E12").Value = "'1.52%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'0.24%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001259"
$ws.Range("E14").Value = "'0.71%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005886"
$ws.Range("E15").Value = "'1.33%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007498"
$ws.Range("E16").Value = "'2,392.15%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.355"
$ws.Range("E17").Value = "'-0.16%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.321"
$ws.Range("E18").Value = "'1.65%"
$ws.Range("D19").Value = "'0.3335"
$ws.Range("D20").Value = "'7.194"
$ws.Range("E20").Value = "'9.50%"
$ws.Range("E21").Value = "'1.70%"
$ws.Range("D22").Value = "'0.2892"
$ws.Range("E22").Value = "'7.87%"
$ws.Range("D23").Value = "'0.04072"
$ws.Range("E23").Value = "'-2.14%"
$ws.Range("D24").Value = "'0.001258"
$ws.Range("E24").Value = "'2.79%"
$ws.Range("D25").Value = "'0.004153"
$ws.Range("E25").Value = "'1.82%"
$ws.Range("D26").Value = "'0.0001273"
$ws.Range("E26").Value = "'-2.06%"
$ws.Range("D38").Value = "'0.02462"
$ws.Range("E38").Value = "'2.23%"
$ws.Range("D39").Value = "'0.05305"
$ws.Range("E39").Value = "'2.01%"
$ws.Range("D40").Value = "'0.007850"
$ws.Range("E40").Value = "'1.24%"
$ws.Range("D41").Value = "'0.1316"
$ws.Range("E41").Value = "'1.38%"
$ws.Range("D42").Value = "'0.006820"
$ws.Range("E42").Value = "'-1.93%"
$ws.Range("D43").Value = "'0.001915"
$ws.Range("E43").Value = "'-0.29%"
$ws.Range("D45").Value = "'0.3061"
$ws.Range("E45").Value = "'0.16%"
$ws.Range("D46").Value = "'0.00006663"
$ws.Range("E46").Value = "'-0.92%"
$ws.Range("E47").Value = "'0.25%"
$ws.Range("D48").Value = "'0.2055"
$ws.Range("E48").Value = "'2,109.84%"
$ws.Range("E49").Value = "'-2.39%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.25%"
